# Actualiza los estadísticos (filas 2-5) de las 3 hojas de parciales con el
# "Nuevo formato 15 jun 2021": más aprobados / menos reprobados y un
# promedio (columna I) ligeramente distinto por hoja.
#
# Para cada hoja, fila: D = Totales (no cambia), E = Aprobados, F = Reprobados,
# G = % Aprobados = ROUND(E/D*100, 2), H = % Reprobados = ROUND(F/D*100, 2),
# I = Promedio, J = Blancos (copia de F), K = % Blancos (copia de H).

$wb = $excel.ActiveWorkbook

# Aprobados / Reprobados (iguales en las 3 hojas) por fila, y el promedio (I)
# que difiere por hoja.
$rows = @(2, 3, 4, 5)
$aprobados  = @{ 2 = 20; 3 = 20; 4 = 35; 5 = 19 }
$reprobados = @{ 2 = 1;  3 = 1;  4 = 4;  5 = 2 }

$promedios = @{
    1 = @{ 2 = 8.7; 3 = 8.6; 4 = 7.6; 5 = 8.7 }   # "1er Parcial"
    2 = @{ 2 = 7.4; 3 = 7.7; 4 = 7.5; 5 = 8.8 }   # "2o Parcial"
    3 = @{ 2 = 8.0; 3 = 8.2; 4 = 7.5; 5 = 8.9 }   # "3er Parcial"
}

for ($sheetIdx = 1; $sheetIdx -le 3; $sheetIdx++) {
    $ws = $wb.Worksheets.Item($sheetIdx)

    foreach ($r in $rows) {
        $total = $ws.Cells.Item($r, 4).Value()   # columna D: Totales

        $aprob  = $aprobados[$r]
        $repro  = $reprobados[$r]
        $pctAprob = [Math]::Round(($aprob / $total) * 100, 2)
        $pctRepro = [Math]::Round(($repro / $total) * 100, 2)
        $promedio = $promedios[$sheetIdx][$r]

        $ws.Cells.Item($r, 5).Value  = $aprob      # E: Aprobados
        $ws.Cells.Item($r, 6).Value  = $repro      # F: Reprobados
        $ws.Cells.Item($r, 7).Value  = $pctAprob   # G: Por_Apro
        $ws.Cells.Item($r, 8).Value  = $pctRepro   # H: Por_Repro
        $ws.Cells.Item($r, 9).Value  = $promedio   # I: Promedio
        $ws.Cells.Item($r, 10).Value = $repro      # J: Blancos
        $ws.Cells.Item($r, 11).Value = $pctRepro   # K: Por_Blan
    }
}
